$d = $word.ActiveDocument

# Locate the "KEY ACHIEVEMENTS AND IMPACT" section bullet paragraphs.
# (There is a duplicate, differently-worded, "Discovered systematic..." line
#  in PROFESSIONAL EXPERIENCE, and a byte-for-byte duplicate of the FEC
#  bullet there too, so we must scope edits to the paragraphs that fall
#  between the "KEY ACHIEVEMENTS AND IMPACT" heading and the next heading
#  rather than doing a document-wide Find/Replace on those lines.)
$count = $d.Paragraphs.Count
$bulletIndexes = New-Object System.Collections.ArrayList
$inSection = $false
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*KEY ACHIEVEMENTS AND IMPACT*") {
        $inSection = $true
        continue
    }
    if ($inSection -and $t -like "*TECHNICAL SKILLS*") {
        break
    }
    if ($inSection -and $t -like "•*") {
        [void]$bulletIndexes.Add($i)
    }
}

Write-Output ("bullet paragraphs: " + ($bulletIndexes -join ","))

# bulletIndexes[0..5] correspond (in order) to the six original bullets:
#   0: Discovered systematic race coding errors...
#   1: Trigonometric algorithm for boundary estimation...
#   2: Built redistricting platform used by thousands...
#   3: Achieved 87% prediction accuracy...
#   4: Built real-time FEC analysis systems... (removed entirely)
#   5: Provided expert testimony and press briefings...

$p1 = $d.Paragraphs.Item($bulletIndexes[0])
$p1.Range.Find.Execute(
    "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions",
    2) | Out-Null

$p2 = $d.Paragraphs.Item($bulletIndexes[1])
$p2.Range.Find.Execute(
    "Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "178% accuracy improvement in racial classification algorithms",
    2) | Out-Null

$p3 = $d.Paragraphs.Item($bulletIndexes[2])
$p3.Range.Find.Execute(
    "Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%",
    2) | Out-Null

$p4 = $d.Paragraphs.Item($bulletIndexes[3])
$p4.Range.Find.Execute(
    "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "`$4.7M savings enabled nonprofit access",
    2) | Out-Null

$p6 = $d.Paragraphs.Item($bulletIndexes[5])
$p6.Range.Find.Execute(
    "Provided expert testimony and press briefings on electoral data integrity and demographic modeling accuracy",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations",
    2) | Out-Null

# Remove the FEC bullet paragraph entirely (paragraph 5 of the section).
$p5 = $d.Paragraphs.Item($bulletIndexes[4])
$p5.Range.Delete()

Write-Output "done"
